# Start evaluation and cleaning of documents
# Clear the leftover analysis columns (C:G) in row 48, which removed two
# now-unused shared strings ("general normative statement" and the
# "Prescription on the need to take on action..." note).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C48:G48").ClearContents()

# Update the saved view/selection to match where the author was working.
$ws.Application.ActiveWindow.ScrollRow = 44
$ws.Range("C48:H48").Select()
